$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 10055.272
$ws.Range("I33").Value = 14419.066
$ws.Range("J33").Value = 704.2857
$ws.Range("K33").Value = 14419.066
$ws.Range("L33").Value = 704.2857
$ws.Range("M33").Value = -14190.066
$ws.Range("N33").Value = -1162.2857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 7999.3335
$ws.Range("I47").Value = 7999.3335
$ws.Range("K47").Value = 7999.3335
$ws.Range("M47").Value = -7027.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 12753.459
$ws.Range("I137").Value = 12706.333
$ws.Range("J137").Value = 12955.429
$ws.Range("K137").Value = 38118.999
$ws.Range("L137").Value = 38866.287
$ws.Range("M137").Value = -35568.999
$ws.Range("N137").Value = -43966.287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 23258294
$ws.Range("J138").Value = 52635944
$ws.Range("L138").Value = 157907832
$ws.Range("N138").Value = -157918112

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4045.1035
$ws.Range("I32").Value = 4019.9285
$ws.Range("K32").Value = 4019.9285
$ws.Range("M32").Value = -3732.9285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3774.2593
$ws.Range("I61").Value = 3440.353
$ws.Range("K61").Value = 3440.353
$ws.Range("M61").Value = -3228.353

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 35523.715
$ws.Range("I74").Value = 43675.82
$ws.Range("J74").Value = 2915.2856
$ws.Range("K74").Value = 43675.82
$ws.Range("L74").Value = 2915.2856
$ws.Range("M74").Value = -42801.82
$ws.Range("N74").Value = -4663.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 35523.715
$ws.Range("I77").Value = 43675.82
$ws.Range("J77").Value = 2915.2856
$ws.Range("K77").Value = 218379.1
$ws.Range("L77").Value = 14576.428
$ws.Range("M77").Value = -214011.1
$ws.Range("N77").Value = -23312.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 24666.666
$ws.Range("J92").Value = 24666.666
$ws.Range("L92").Value = 24666.666
$ws.Range("N92").Value = -29658.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 26736.143
$ws.Range("I132").Value = 1727.8113
$ws.Range("J132").Value = 468550
$ws.Range("K132").Value = 5183.4339
$ws.Range("L132").Value = 1405650
$ws.Range("M132").Value = -2653.4339
$ws.Range("N132").Value = -1410710

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3774.2593
$ws.Range("I136").Value = 3440.353
$ws.Range("K136").Value = 10321.059
$ws.Range("M136").Value = -7771.059000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 198000
$ws.Range("J138").Value = 198000
$ws.Range("L138").Value = 198000
$ws.Range("N138").Value = -208280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 199990
$ws.Range("J139").Value = 199990
$ws.Range("L139").Value = 199990
$ws.Range("N139").Value = -210270

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 43065
$ws.Range("J58").Value = 48542
$ws.Range("L58").Value = 48542
$ws.Range("N58").Value = -49130

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4756.375
$ws.Range("I99").Value = 4430
$ws.Range("J99").Value = 4803
$ws.Range("K99").Value = 4430
$ws.Range("L99").Value = 4803
$ws.Range("M99").Value = -2932
$ws.Range("N99").Value = -7799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 86929.664
$ws.Range("J52").Value = 86929.664
$ws.Range("L52").Value = 86929.664
$ws.Range("N52").Value = -87517.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3860.5789
$ws.Range("I62").Value = 3883.625
$ws.Range("J62").Value = 3843.818
$ws.Range("K62").Value = 3883.625
$ws.Range("L62").Value = 3843.818
$ws.Range("M62").Value = -3259.625
$ws.Range("N62").Value = -5091.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3860.5789
$ws.Range("I65").Value = 3883.625
$ws.Range("J65").Value = 3843.818
$ws.Range("K65").Value = 19418.125
$ws.Range("L65").Value = 19219.09
$ws.Range("M65").Value = -16298.125
$ws.Range("N65").Value = -25459.09

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4722.4375
$ws.Range("I99").Value = 4326.0713
$ws.Range("J99").Value = 7497
$ws.Range("K99").Value = 4326.0713
$ws.Range("L99").Value = 7497
$ws.Range("M99").Value = -2828.0713
$ws.Range("N99").Value = -10493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H115").Value = 66996.336
$ws.Range("J115").Value = 69990
$ws.Range("L115").Value = 69990
$ws.Range("N115").Value = -72340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4722.4375
$ws.Range("I126").Value = 4326.0713
$ws.Range("J126").Value = 7497
$ws.Range("K126").Value = 12978.2139
$ws.Range("L126").Value = 22491
$ws.Range("M126").Value = -10508.2139
$ws.Range("N126").Value = -27431

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3682.3242
$ws.Range("I132").Value = 3448.2727
$ws.Range("K132").Value = 10344.8181
$ws.Range("M132").Value = -7814.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 995.2258
$ws.Range("J131").Value = 1204.4706
$ws.Range("L131").Value = 3613.4118
$ws.Range("N131").Value = -13693.4118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4010.4443
$ws.Range("J102").Value = 7313.923
$ws.Range("L102").Value = 7313.923
$ws.Range("N102").Value = -10557.923

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3999
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1182.8
$ws.Range("I132").Value = 1150
$ws.Range("J132").Value = 1204.6666
$ws.Range("K132").Value = 3450
$ws.Range("L132").Value = 3613.9998
$ws.Range("M132").Value = -920
$ws.Range("N132").Value = -8673.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4399.778
$ws.Range("I7").Value = 3399.7693
$ws.Range("K7").Value = 3399.7693
$ws.Range("M7").Value = -3287.7693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1083.9231
$ws.Range("I16").Value = 1007.5833
$ws.Range("K16").Value = 1007.5833
$ws.Range("M16").Value = -837.5833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4399.778
$ws.Range("I126").Value = 3399.7693
$ws.Range("K126").Value = 10199.3079
$ws.Range("M126").Value = -7729.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 170000
$ws.Range("J46").Value = 170000
$ws.Range("L46").Value = 170000
$ws.Range("N46").Value = -170462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7283.5713
$ws.Range("I62").Value = 5497
$ws.Range("K62").Value = 5497
$ws.Range("M62").Value = -4873

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7283.5713
$ws.Range("I65").Value = 5497
$ws.Range("K65").Value = 27485
$ws.Range("M65").Value = -24365

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 63143
$ws.Range("J124").Value = 63143
$ws.Range("L124").Value = 63143
$ws.Range("N124").Value = -72963

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 70938.8
$ws.Range("J125").Value = 70938.8
$ws.Range("L125").Value = 70938.8
$ws.Range("N125").Value = -80778.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2296.2856
$ws.Range("I126").Value = 2018.8
$ws.Range("K126").Value = 6056.4
$ws.Range("M126").Value = -3586.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 170000
$ws.Range("J134").Value = 170000
$ws.Range("L134").Value = 510000
$ws.Range("N134").Value = -515070

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3058.5781
$ws.Range("I136").Value = 1969.234
$ws.Range("J136").Value = 6070.294
$ws.Range("K136").Value = 5907.701999999999
$ws.Range("L136").Value = 18210.882
$ws.Range("M136").Value = -3357.701999999999
$ws.Range("N136").Value = -23310.882

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 121715
$ws.Range("J137").Value = 121715
$ws.Range("L137").Value = 121715
$ws.Range("N137").Value = -131915
